$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81; Excel shifts rows 81-127 down to 82-128
# and preserves formatting (e.g. the date style on column D).
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with the new weekly record.
$ws.Range("A81").Value = 4
$ws.Range("B81").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C81").Value = "Los Lagos"
$ws.Range("D81").Value = 44452
$ws.Range("E81").Value = 10
$ws.Range("F81").Value = 100112017
$ws.Range("G81").Value = "Apio"
$ws.Range("H81").Value = "Americana (o)"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 25
$ws.Range("K81").Value = 12000
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = 12000
$ws.Range("N81").Value = "$/docena de matas"
$ws.Range("O81").Value = "Región de Coquimbo"
$ws.Range("P81").Value = 2000
$ws.Range("Q81").Value = 6
$ws.Range("R81").Value = "Hortaliza"
